$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellAddr, $Text)
    $ws.Range("Z1").NumberFormat = "@"
    $ws.Range("Z1").Value = $Text
    $ws.Range("Z1").Copy() | Out-Null
    $ws.Range($CellAddr).PasteSpecial(-4163) | Out-Null
    $ws.Range("Z1").Clear() | Out-Null
}

# Direct text/string updates (safe as text: not single-dot numeric-looking)
$ws.Range("D2").Value = '51.762.77'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '2.778.47'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("E7").Value = '  -1.60%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("E10").Value = '  -2.21%  '
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("E13").Value = '  -1.81%  '
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").Value = '3.213.18'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '2.779.07'
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = '51.687.05'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").Value = '0.0₃0969'
$ws.Range("E22").Value = '  -2.15%  '
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  +17.94%  '
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  +1.19%  '
$ws.Range("E31").Value = '  +5.63%  '
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("E34").Value = '  -8.71%  '
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("E36").Value = '  -6.72%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  +1.96%  '
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("E40").Value = '  -3.46%  '
$ws.Range("E41").Value = '  +2.54%  '
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  -5.96%  '
$ws.Range("E45").Value = '  -6.16%  '
$ws.Range("D46").Value = '2.083.58'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("E50").Value = '  -5.88%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E51").Value = '  +0.69%  '

# Updates that must be forced to stay text (values that look like plain numbers)
Set-TextValue "D5" '356.35'
Set-TextValue "D6" '108.90'
Set-TextValue "D8" '1.00'
Set-TextValue "D9" '0.585'
Set-TextValue "D10" '39.63'
Set-TextValue "D12" '0.0845'
Set-TextValue "D13" '19.44'
Set-TextValue "D14" '7.60'
Set-TextValue "D19" '7.44'
Set-TextValue "D20" '3.09'
Set-TextValue "D21" '13.12'
Set-TextValue "D23" '70.19'
Set-TextValue "D24" '268.71'
Set-TextValue "D26" '26.37'
Set-TextValue "D28" '0.165'
Set-TextValue "D29" '10.19'
Set-TextValue "D30" '2.28'
Set-TextValue "D31" '6.23'
Set-TextValue "D32" '34.73'
Set-TextValue "D34" '0.0451'
Set-TextValue "D35" '0.0839'
Set-TextValue "D37" '1.00'
Set-TextValue "D38" '18.62'
Set-TextValue "D39" '3.13'
Set-TextValue "D40" '1.95'
Set-TextValue "D43" '2.21'
Set-TextValue "D44" '119.23'
Set-TextValue "D45" '21.71'
Set-TextValue "D47" '3.27'
Set-TextValue "D49" '0.942'
Set-TextValue "D50" '5.57'
Set-TextValue "D51" '0.191'

$excel.CutCopyMode = 0